$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set cell values (cell style/format is preserved automatically)
$ws.Range("D4").Value = 5
$ws.Range("F4").Value = 5
$ws.Range("G4").Value = 5
$ws.Range("H4").Value = 5

$ws.Range("D8").Value = 5
$ws.Range("F8").Value = 5
$ws.Range("H8").Value = 5

$ws.Range("D10").Value = 5
$ws.Range("F10").Value = 5

$ws.Range("D13").Value = 5

$ws.Range("D15").Value = 5
$ws.Range("F15").Value = 5

$ws.Range("D16").Value = 5
$ws.Range("F16").Value = 5

$ws.Range("D19").Value = 5
$ws.Range("F19").Value = 5
$ws.Range("G19").Value = 5

$ws.Range("D22").Value = 5
$ws.Range("F22").Value = 5
$ws.Range("G22").Value = 5

$ws.Range("D26").Value = 5

# Update the active selection to D5
$ws.Range("D5").Select()
